$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("P3").Value = 311
$ws.Range("Q3").Value = 0
$ws.Range("S3").Value = 86
$ws.Range("T3").Value = 86
$ws.Range("U3").Value = 0
$ws.Rows.Item(3).Hidden = $true

# Row 4
$ws.Range("L4").Value = -1
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("P4").Value = 319
$ws.Range("Q4").Value = 0
$ws.Range("S4").Value = 4
$ws.Range("T4").Value = 4
$ws.Range("U4").Value = 0
$ws.Rows.Item(4).Hidden = $true

# Row 8
$ws.Range("M8").Value = 21.25
$ws.Range("N8").Value = 12.75
$ws.Range("P8").Value = 1
$ws.Range("Q8").Value = 1
$ws.Range("U8").Value = 1

# Row 10
$ws.Range("L10").Value = -1
$ws.Range("M10").Value = 14.25
$ws.Range("N10").Value = 8.550000000000001
$ws.Range("S10").Value = 3
$ws.Range("T10").Value = 3
$ws.Range("U10").Value = 6

# Row 11
$ws.Range("M11").Value = 0
$ws.Range("N11").Value = 0
$ws.Range("P11").Value = 32
$ws.Range("Q11").Value = 0
$ws.Range("S11").Value = 1
$ws.Range("T11").Value = 1
$ws.Range("U11").Value = 0
$ws.Rows.Item(11).Hidden = $true

# Row 12
$ws.Range("L12").Value = 2
$ws.Range("M12").Value = 49.6
$ws.Range("N12").Value = 29.76
$ws.Range("S12").Value = 4
$ws.Range("T12").Value = 4
$ws.Range("U12").Value = 16

# Row 13
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("P13").Value = 32
$ws.Range("Q13").Value = 0
$ws.Range("U13").Value = 0
$ws.Rows.Item(13).Hidden = $true

# Row 15
$ws.Range("L15").Value = -1
$ws.Range("M15").Value = 0
$ws.Range("N15").Value = 0
$ws.Range("P15").Value = 214
$ws.Range("Q15").Value = 0
$ws.Range("S15").Value = 1
$ws.Range("T15").Value = 1
$ws.Range("U15").Value = 0
$ws.Rows.Item(15).Hidden = $true

# Row 16
$ws.Range("M16").Value = 0
$ws.Range("N16").Value = 0
$ws.Range("P16").Value = 52
$ws.Range("Q16").Value = 0
$ws.Range("S16").Value = 1
$ws.Range("T16").Value = 1
$ws.Range("U16").Value = 0
$ws.Rows.Item(16).Hidden = $true

# Row 17
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("P17").Value = 69
$ws.Range("Q17").Value = 0
$ws.Range("S17").Value = 8
$ws.Range("T17").Value = 8
$ws.Range("U17").Value = 0
$ws.Rows.Item(17).Hidden = $true

# Row 18
$ws.Range("M18").Value = 0
$ws.Range("N18").Value = 0
$ws.Range("P18").Value = 20
$ws.Range("Q18").Value = 0
$ws.Range("S18").Value = 8
$ws.Range("T18").Value = 8
$ws.Range("U18").Value = 0
$ws.Rows.Item(18).Hidden = $true

# Row 19
$ws.Range("L19").Value = 1
$ws.Range("M19").Value = 28.8
$ws.Range("N19").Value = 17.28
$ws.Range("P19").Value = 1
$ws.Range("Q19").Value = 5
$ws.Range("S19").Value = 11
$ws.Range("T19").Value = 11
$ws.Range("U19").Value = 16

# Row 20
$ws.Range("L20").Value = 1
$ws.Range("M20").Value = 262.5
$ws.Range("N20").Value = 157.5
$ws.Range("P20").Value = 2
$ws.Range("Q20").Value = 2
$ws.Range("S20").Value = 1
$ws.Range("T20").Value = 1
$ws.Range("U20").Value = 3

# Row 21
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = 0
$ws.Range("P21").Value = 60
$ws.Range("Q21").Value = 0
$ws.Range("U21").Value = 0
$ws.Rows.Item(21).Hidden = $true

# Row 22
$ws.Range("L22").Value = -1
$ws.Range("P22").Value = 226
$ws.Range("Q22").Value = 0
$ws.Range("U22").Value = 0
$ws.Rows.Item(22).Hidden = $true

# Row 23
$ws.Range("L23").Value = -1
$ws.Range("P23").Value = 14
$ws.Range("Q23").Value = 0
$ws.Range("S23").Value = 14
$ws.Range("T23").Value = 14

# Row 25
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("P25").Value = 8
$ws.Range("Q25").Value = 0
$ws.Range("U25").Value = 0
$ws.Rows.Item(25).Hidden = $true

# Row 27
$ws.Range("L27").Value = -1
$ws.Range("M27").Value = 116.2
$ws.Range("N27").Value = 69.72
$ws.Range("S27").Value = 4
$ws.Range("T27").Value = 4
$ws.Range("U27").Value = 8

# Row 28
$ws.Range("L28").Value = -1

# Row 29
$ws.Range("P29").Value = 16
$ws.Range("Q29").Value = 0
$ws.Range("S29").Value = 4
$ws.Range("T29").Value = 4
$ws.Range("U29").Value = 0
$ws.Rows.Item(29).Hidden = $true

# Row 30
$ws.Range("L30").Value = -5

# Row 31
$ws.Range("L31").Value = 1
$ws.Range("M31").Value = 0
$ws.Range("N31").Value = 0
$ws.Range("P31").Value = 223
$ws.Range("Q31").Value = 0
$ws.Range("S31").Value = 2
$ws.Range("T31").Value = 2
$ws.Range("U31").Value = 0
$ws.Rows.Item(31).Hidden = $true

# Row 32
$ws.Range("P32").Value = 84
$ws.Range("Q32").Value = 0
$ws.Range("S32").Value = 1
$ws.Range("T32").Value = 1
$ws.Range("U32").Value = 0
$ws.Rows.Item(32).Hidden = $true

# Row 33
$ws.Range("M33").Value = 0
$ws.Range("N33").Value = 0
$ws.Range("P33").Value = 69
$ws.Range("Q33").Value = 0
$ws.Range("S33").Value = 2
$ws.Range("T33").Value = 2
$ws.Range("U33").Value = 0
$ws.Rows.Item(33).Hidden = $true

# Row 34
$ws.Range("L34").Value = -7
$ws.Range("P34").Value = 35
$ws.Range("Q34").Value = 0
$ws.Range("U34").Value = 0
$ws.Rows.Item(34).Hidden = $true

# Row 35
$ws.Range("P35").Value = 5
$ws.Range("Q35").Value = 0
$ws.Range("S35").Value = 5
$ws.Range("T35").Value = 5

# Row 40
$ws.Range("C40").Value = 89

# Row 42
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "841.82€"
$ws.Range("C41").Copy() | Out-Null
$ws.Range("C42").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 51
$ws.Range("C51").Value = -13
